# TalDoor BOM update — swap the 1000uF cap (C15) for a shorter/cheaper part,
# note the old part number, and highlight the Qty column for the populated rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 corresponds to ref "C15," (1000uF polarised capacitor).
# The old, taller part is being swapped for a lower-profile one.
$ws.Range("F11").Value = "1572-1665-ND"
$ws.Range("G11").Value = 0.283
$ws.Range("I11").Value = "OLD taller one: 493-5907-1-ND"

# Highlight the Qty column (B) for the populated BOM rows with a light
# green fill + centered text, skipping the rows that were left as-is
# (B16, B19, B31, B32).
$qtyRanges = @("B5:B15", "B17:B18", "B20:B30", "B33:B34")
foreach ($addr in $qtyRanges) {
    $rng = $ws.Range($addr)
    $rng.HorizontalAlignment = -4108  # xlCenter
    $rng.Interior.ThemeColor = 8      # xlThemeColorAccent6
    $rng.Interior.TintAndShade = 0.6
}

$ws.Range("I11").Select()
